$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value (serial 45179 = 2023-09-10) for every
# data row (rows 2 through 215). Update all of them by one day to 45180 (2023-09-11).
$ws.Range("C2:C215").Value = 45180
